$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: comment text corrected ---
$ws.Range("H23").Value = "could not calibrate the eyes "

# --- Rows 27-30: fill in dates, gender, comments (order/G already populated) ---
$ws.Range("B27").Value = 42942
$ws.Range("C27").Value = 41153
$ws.Range("D27").Value = "F"
$ws.Range("H27").Value = "session went well"

$ws.Range("B28").Value = 42942
$ws.Range("C28").Value = 41779
$ws.Range("D28").Value = "M"
$ws.Range("H28").Value = "session went well"

$ws.Range("B29").Value = 42942
$ws.Range("C29").Value = 41856
$ws.Range("D29").Value = "M"
$ws.Range("H29").Value = "session went well"

$ws.Range("B30").Value = 42942
$ws.Range("C30").Value = 41608
$ws.Range("D30").Value = "M"
$ws.Range("H30").Value = "played with headphones during a couple of trials, and stood up for a couple of trials, but majority of session went well "

# --- Rows 31-33: same as above, plus B column font needs to normalize to Arial 10 (matches rows 2-30) ---
$ws.Range("B31").Value = 42942
$ws.Range("B31").Font.Name = "Arial"
$ws.Range("B31").Font.Size = 10
$ws.Range("C31").Value = 41552
$ws.Range("D31").Value = "M"
$ws.Range("H31").Value = "could not calibrate the eyes "

$ws.Range("B32").Value = 42942
$ws.Range("B32").Font.Name = "Arial"
$ws.Range("B32").Font.Size = 10
$ws.Range("C32").Value = 40966
$ws.Range("D32").Value = "F"
$ws.Range("H32").Value = "session went well"

$ws.Range("B33").Value = 42942
$ws.Range("B33").Font.Name = "Arial"
$ws.Range("B33").Font.Size = 10
$ws.Range("C33").Value = 41191
$ws.Range("D33").Value = "M"
$ws.Range("H33").Value = "session went well"

# --- Rows 34-36: new run dates, birthdays, gender, comments ---
$ws.Range("B34").Value = 42944
$ws.Range("C34").Value = 41307
$ws.Range("D34").Value = "F"
$ws.Range("H34").Value = "session went well"

$ws.Range("B35").Value = 42944
$ws.Range("C35").Value = 40928
$ws.Range("D35").Value = "F"
$ws.Range("H35").Value = "session went well"

$ws.Range("B36").Value = 42944
$ws.Range("C36").Value = 40834
$ws.Range("D36").Value = "F"
$ws.Range("H36").Value = "session went well"

# --- Selection / active cell moved to H38 (scrolled down a row) ---
$ws.Range("H38").Select()
